$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Add the new "metadata" worksheet, positioned right after "data" ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$newSheet.Name = "metadata"

# Reuse the existing header style (bold + border + center/top alignment)
# from the "data" sheet's header row instead of building a brand new style.
$dataSheet.Range("B1:F1").Copy($newSheet.Range("B1:F1"))
$dataSheet.Range("F1").Copy($newSheet.Range("G1"))
$dataSheet.Range("A2").Copy($newSheet.Range("A2"))

# Header row
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Diabetes - neonatal onset"
$newSheet.Range("C2").Value = 293
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "2.34"
$newSheet.Range("E2").Value = "2021-05-19T14:06:19.357844Z"
$newSheet.Range("F2").Value = "2021-10-05 14:19:58.675969"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/293/?format=json"

# --- Refresh the "time_taken" timestamps on the "data" sheet (col F, rows 2-35) ---
$dataSheet.Range("F2").Value = "2021-10-05 14:19:58.679594"
$dataSheet.Range("F3").Value = "2021-10-05 14:19:58.679602"
$dataSheet.Range("F4").Value = "2021-10-05 14:19:58.679605"
$dataSheet.Range("F5").Value = "2021-10-05 14:19:58.679608"
$dataSheet.Range("F6").Value = "2021-10-05 14:19:58.679611"
$dataSheet.Range("F7").Value = "2021-10-05 14:19:58.679614"
$dataSheet.Range("F8").Value = "2021-10-05 14:19:58.679616"
$dataSheet.Range("F9").Value = "2021-10-05 14:19:58.679619"
$dataSheet.Range("F10").Value = "2021-10-05 14:19:58.679621"
$dataSheet.Range("F11").Value = "2021-10-05 14:19:58.679624"
$dataSheet.Range("F12").Value = "2021-10-05 14:19:58.679626"
$dataSheet.Range("F13").Value = "2021-10-05 14:19:58.679629"
$dataSheet.Range("F14").Value = "2021-10-05 14:19:58.679631"
$dataSheet.Range("F15").Value = "2021-10-05 14:19:58.679634"
$dataSheet.Range("F16").Value = "2021-10-05 14:19:58.679636"
$dataSheet.Range("F17").Value = "2021-10-05 14:19:58.679639"
$dataSheet.Range("F18").Value = "2021-10-05 14:19:58.679642"
$dataSheet.Range("F19").Value = "2021-10-05 14:19:58.679644"
$dataSheet.Range("F20").Value = "2021-10-05 14:19:58.679647"
$dataSheet.Range("F21").Value = "2021-10-05 14:19:58.679650"
$dataSheet.Range("F22").Value = "2021-10-05 14:19:58.679652"
$dataSheet.Range("F23").Value = "2021-10-05 14:19:58.679655"
$dataSheet.Range("F24").Value = "2021-10-05 14:19:58.679658"
$dataSheet.Range("F25").Value = "2021-10-05 14:19:58.679660"
$dataSheet.Range("F26").Value = "2021-10-05 14:19:58.679663"
$dataSheet.Range("F27").Value = "2021-10-05 14:19:58.679665"
$dataSheet.Range("F28").Value = "2021-10-05 14:19:58.679668"
$dataSheet.Range("F29").Value = "2021-10-05 14:19:58.679670"
$dataSheet.Range("F30").Value = "2021-10-05 14:19:58.679673"
$dataSheet.Range("F31").Value = "2021-10-05 14:19:58.679676"
$dataSheet.Range("F32").Value = "2021-10-05 14:19:58.679678"
$dataSheet.Range("F33").Value = "2021-10-05 14:19:58.679681"
$dataSheet.Range("F34").Value = "2021-10-05 14:19:58.679683"
$dataSheet.Range("F35").Value = "2021-10-05 14:19:58.679686"

# Re-select the "data" sheet as active, matching the original workbook view
$dataSheet.Activate()

Write-Host "metadata sheet added and timestamps refreshed"
